$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells in this sheet are stored as text (inline strings),
# even when they look numeric (e.g. "540.53") or use "." as a thousands
# separator (e.g. "58.961.93"). Force text format so Excel COM does not
# silently coerce these into floating point numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.961.93'
$ws.Range('E2').Value = '  -3.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.228.46'
$ws.Range('E3').Value = '  -3.89%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.53'
$ws.Range('E5').Value = '  -4.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.40'
$ws.Range('E6').Value = '  -8.13%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.226.06'
$ws.Range('E8').Value = '  -3.95%  '
$ws.Range('E9').Value = '  -4.17%  '
$ws.Range('E10').Value = '  -3.94%  '
$ws.Range('E12').Value = '  -4.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.781.87'
$ws.Range('E13').Value = '  -3.90%  '
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.01'
$ws.Range('E15').Value = '  -6.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.230.98'
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('E17').Value = '  -5.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.003.47'
$ws.Range('E19').Value = '  -6.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.39'
$ws.Range('E20').Value = '  -5.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.23'
$ws.Range('E21').Value = '  -6.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '363.35'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.61'
$ws.Range('E24').Value = '  -6.16%  '
$ws.Range('E25').Value = '  -6.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.361.78'
$ws.Range('E26').Value = '  -3.97%  '
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0972'
$ws.Range('E28').Value = '  -9.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.16'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.93'
$ws.Range('E32').Value = '  -6.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.11'
$ws.Range('E33').Value = '  -7.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.93'
$ws.Range('E34').Value = '  -3.92%  '
$ws.Range('E35').Value = '  -6.10%  '
$ws.Range('E36').Value = '  -7.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.32'
$ws.Range('E37').Value = '  -5.23%  '
$ws.Range('E38').Value = '  -4.92%  '
$ws.Range('E39').Value = '  -6.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.36'
$ws.Range('E40').Value = '  -8.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0710'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.257.50'
$ws.Range('E42').Value = '  -4.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.14'
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('E44').Value = '  -5.62%  '
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('E46').Value = '  -5.66%  '
$ws.Range('E47').Value = '  -6.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.303.57'
$ws.Range('E49').Value = '  -7.37%  '
$ws.Range('E50').Value = '  -4.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.93'
$ws.Range('E51').Value = '  -6.67%  '
